$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rule row for Complaints: deny access to non-participants when the
# "restricted" flag is set on the complaint.
$ws.Range("B28").Value = "Complaint – Restricted Flag"
$ws.Range("C28").Value = "COMPLAINT"
$ws.Range("D28").Value = "restricted"
$ws.Range("G28").Value = "deny read to *"

# New rule row for Case Files: deny access to non-participants when the
# "restricted" flag is set on the case file.
$ws.Range("B29").Value = "Case File – Restricted Flag"
$ws.Range("C29").Value = "CASE_FILE"
$ws.Range("D29").Value = "restricted"
$ws.Range("G29").Value = "deny read to *"

$ws.Range("B30").Select()
